$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1863799283154122
$ws.Range("C2").Value = 0.5663082437275986
$ws.Range("J2").Value = 0.01075268817204301
$ws.Range("P2").Value = 0.1362007168458781
$ws.Range("S2").Value = 0.1003584229390681

$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.03067484662576687
$ws.Range("J3").Value = 0.03680981595092025
$ws.Range("P3").Value = 0.8343558282208589
$ws.Range("S3").Value = 0.09202453987730061

$ws.Range("J4").Value = 0.1142857142857143
$ws.Range("P4").Value = 0.6285714285714286
$ws.Range("S4").Value = 0.2571428571428571

$ws.Range("B6").Value = 0.05803571428571429
$ws.Range("D6").Value = 0.01339285714285714
$ws.Range("F6").Value = 0.1071428571428571
$ws.Range("J6").Value = 0.1785714285714286
$ws.Range("O6").Value = 0.01785714285714286
$ws.Range("Q6").Value = 0.1785714285714286
$ws.Range("R6").Value = 0.0625
$ws.Range("S6").Value = 0.3839285714285715

$ws.Range("B7").Value = 0.1073446327683616
$ws.Range("D7").Value = 0.01129943502824859
$ws.Range("E7").Value = 0.005649717514124294
$ws.Range("F7").Value = 0.04519774011299435
$ws.Range("J7").Value = 0.1186440677966102
$ws.Range("O7").Value = 0.005649717514124294
$ws.Range("Q7").Value = 0.1525423728813559
$ws.Range("R7").Value = 0.1638418079096045
$ws.Range("S7").Value = 0.3898305084745763

$ws.Range("B8").Value = 0.1020881670533643
$ws.Range("D8").Value = 0.009280742459396751
$ws.Range("E8").Value = 0.002320185614849188
$ws.Range("F8").Value = 0.07192575406032482
$ws.Range("J8").Value = 0.06728538283062645
$ws.Range("O8").Value = 0.01624129930394431
$ws.Range("Q8").Value = 0.1624129930394431
$ws.Range("R8").Value = 0.1554524361948956
$ws.Range("S8").Value = 0.4129930394431555

$ws.Range("B9").Value = 0.08300395256916997
$ws.Range("D9").Value = 0.01185770750988142
$ws.Range("E9").Value = 0.01185770750988142
$ws.Range("F9").Value = 0.06719367588932806
$ws.Range("J9").Value = 0.1225296442687747
$ws.Range("O9").Value = 0.02371541501976284
$ws.Range("Q9").Value = 0.150197628458498
$ws.Range("R9").Value = 0.1146245059288538
$ws.Range("S9").Value = 0.4150197628458498

$ws.Range("B10").Value = 0.09709480122324159
$ws.Range("D10").Value = 0.01758409785932722
$ws.Range("E10").Value = 0.001529051987767584
$ws.Range("F10").Value = 0.06651376146788991
$ws.Range("J10").Value = 0.1131498470948012
$ws.Range("O10").Value = 0.01452599388379205
$ws.Range("Q10").Value = 0.2010703363914373
$ws.Range("R10").Value = 0.1230886850152905
$ws.Range("S10").Value = 0.3654434250764526

$ws.Range("G11").Value = 0.1779359430604982
$ws.Range("J11").Value = 0.09608540925266904
$ws.Range("K11").Value = 0.2384341637010676
$ws.Range("L11").Value = 0.4697508896797153
$ws.Range("S11").Value = 0.01779359430604982

$ws.Range("G12").Value = 0.7238805970149254
$ws.Range("J12").Value = 0.2313432835820896
$ws.Range("S12").Value = 0.04477611940298507

$ws.Range("G13").Value = 0.6785714285714286
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.03571428571428571

$ws.Range("F15").Value = 0.01265822784810127
$ws.Range("H15").Value = 0.2151898734177215
$ws.Range("I15").Value = 0.04641350210970464
$ws.Range("J15").Value = 0.350210970464135
$ws.Range("K15").Value = 0.08016877637130802
$ws.Range("M15").Value = 0.008438818565400843
$ws.Range("O15").Value = 0.09282700421940929
$ws.Range("S15").Value = 0.1940928270042194

$ws.Range("H16").Value = 0.1711229946524064
$ws.Range("I16").Value = 0.1283422459893048
$ws.Range("J16").Value = 0.3422459893048128
$ws.Range("K16").Value = 0.09090909090909091
$ws.Range("M16").Value = 0.0374331550802139
$ws.Range("O16").Value = 0.106951871657754
$ws.Range("S16").Value = 0.1229946524064171

$ws.Range("F17").Value = 0.01157407407407407
$ws.Range("H17").Value = 0.150462962962963
$ws.Range("I17").Value = 0.1388888888888889
$ws.Range("J17").Value = 0.4421296296296297
$ws.Range("K17").Value = 0.07175925925925926
$ws.Range("M17").Value = 0.009259259259259259
$ws.Range("N17").Value = 0.002314814814814815
$ws.Range("O17").Value = 0.05092592592592592
$ws.Range("S17").Value = 0.1226851851851852

$ws.Range("F18").Value = 0.02013422818791946
$ws.Range("H18").Value = 0.1644295302013423
$ws.Range("I18").Value = 0.1006711409395973
$ws.Range("J18").Value = 0.4194630872483222
$ws.Range("K18").Value = 0.08389261744966443
$ws.Range("M18").Value = 0.02013422818791946
$ws.Range("N18").Value = 0.003355704697986577
$ws.Range("O18").Value = 0.06375838926174497
$ws.Range("S18").Value = 0.1241610738255034

$ws.Range("F19").Value = 0.0132398753894081
$ws.Range("H19").Value = 0.1853582554517134
$ws.Range("I19").Value = 0.1004672897196262
$ws.Range("J19").Value = 0.3933021806853583
$ws.Range("K19").Value = 0.09501557632398754
$ws.Range("M19").Value = 0.02725856697819315
$ws.Range("O19").Value = 0.06853582554517133
$ws.Range("S19").Value = 0.1168224299065421
